# Auto-generated edit script: applies updated commodity-price / profit
# figures to each crafting-class sheet per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7576985.5
$ws.Range("J17").Value = 7576985.5
$ws.Range("L17").Value = 22730956.5
$ws.Range("N17").Value = -22731292.5
$ws.Range("H86").Value = 35277.645
$ws.Range("I86").Value = 3740
$ws.Range("J86").Value = 52798.555
$ws.Range("K86").Value = 3740
$ws.Range("L86").Value = 52798.555
$ws.Range("M86").Value = -2617
$ws.Range("N86").Value = -55044.555
$ws.Range("H88").Value = 3360.889
$ws.Range("I88").Value = 3533
$ws.Range("J88").Value = 3274.8333
$ws.Range("K88").Value = 3533
$ws.Range("L88").Value = 3274.8333
$ws.Range("M88").Value = -3127
$ws.Range("N88").Value = -4086.8333
$ws.Range("H89").Value = 35277.645
$ws.Range("I89").Value = 3740
$ws.Range("J89").Value = 52798.555
$ws.Range("K89").Value = 18700
$ws.Range("L89").Value = 263992.775
$ws.Range("M89").Value = -13084
$ws.Range("N89").Value = -275224.775
$ws.Range("H91").Value = 3360.889
$ws.Range("I91").Value = 3533
$ws.Range("J91").Value = 3274.8333
$ws.Range("K91").Value = 3533
$ws.Range("L91").Value = 3274.8333
$ws.Range("M91").Value = -2129
$ws.Range("N91").Value = -6082.8333
$ws.Range("H112").Value = 3959.2778
$ws.Range("I112").Value = 1191.75
$ws.Range("J112").Value = 4750
$ws.Range("K112").Value = 3575.25
$ws.Range("L112").Value = 14250
$ws.Range("M112").Value = -2467.25
$ws.Range("N112").Value = -16466
$ws.Range("H137").Value = 20816.273
$ws.Range("I137").Value = 30197.715
$ws.Range("K137").Value = 90593.145
$ws.Range("M137").Value = -88043.145

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 354182.03
$ws.Range("I122").Value = 3164.5862
$ws.Range("K122").Value = 9493.758600000001
$ws.Range("M122").Value = -7043.758600000001

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 18075.56
$ws.Range("I99").Value = 27854.5
$ws.Range("K99").Value = 27854.5
$ws.Range("M99").Value = -26356.5
$ws.Range("H107").Value = 5499.5
$ws.Range("I107").Value = 10000
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 10000
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = -8080
$ws.Range("N107").Value = -4839
$ws.Range("H140").Value = 76081
$ws.Range("J140").Value = 76081
$ws.Range("L140").Value = 76081
$ws.Range("N140").Value = -86441

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1044
$ws.Range("I16").Value = 1044
$ws.Range("K16").Value = 1044
$ws.Range("M16").Value = -757
$ws.Range("H94").Value = 3068.375
$ws.Range("I94").Value = 1948.5
$ws.Range("J94").Value = 3441.6667
$ws.Range("K94").Value = 1948.5
$ws.Range("L94").Value = 3441.6667
$ws.Range("M94").Value = -1497.5
$ws.Range("N94").Value = -4343.6667
$ws.Range("H107").Value = 9093.666999999999
$ws.Range("I107").Value = 11464.477
$ws.Range("J107").Value = 795.8333
$ws.Range("K107").Value = 11464.477
$ws.Range("L107").Value = 795.8333
$ws.Range("M107").Value = -9544.477000000001
$ws.Range("N107").Value = -4635.8333
$ws.Range("H113").Value = 1044
$ws.Range("I113").Value = 1044
$ws.Range("K113").Value = 1044
$ws.Range("M113").Value = 1126
$ws.Range("H118").Value = 74995
$ws.Range("J118").Value = 74995
$ws.Range("L118").Value = 74995
$ws.Range("N118").Value = -78309

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 108.333336
$ws.Range("I2").Value = 125
$ws.Range("K2").Value = 750
$ws.Range("M2").Value = -637
$ws.Range("H113").Value = 14774.75
$ws.Range("J113").Value = 28249.5
$ws.Range("L113").Value = 84748.5
$ws.Range("N113").Value = -89088.5
$ws.Range("H129").Value = 22224518
$ws.Range("I129").Value = 769.8333
$ws.Range("J129").Value = 37040350
$ws.Range("K129").Value = 2309.4999
$ws.Range("L129").Value = 111121050
$ws.Range("M129").Value = 2690.5001
$ws.Range("N129").Value = -111131050
$ws.Range("H137").Value = 5244.9414
$ws.Range("I137").Value = 1267.75
$ws.Range("J137").Value = 14790.2
$ws.Range("K137").Value = 3803.25
$ws.Range("L137").Value = 44370.60000000001
$ws.Range("M137").Value = 1296.75
$ws.Range("N137").Value = -54570.60000000001

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9790.182000000001
$ws.Range("J80").Value = 2500
$ws.Range("L80").Value = 2500
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 9790.182000000001
$ws.Range("J83").Value = 2500
$ws.Range("L83").Value = 12500
$ws.Range("N83").Value = -22484
$ws.Range("H97").Value = 5757.2085
$ws.Range("I97").Value = 7487.0557
$ws.Range("J97").Value = 567.6667
$ws.Range("K97").Value = 7487.0557
$ws.Range("L97").Value = 567.6667
$ws.Range("M97").Value = -6991.0557
$ws.Range("N97").Value = -1559.6667
$ws.Range("H113").Value = 34811.2
$ws.Range("I113").Value = 34811.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 34811.2
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -32641.2
$ws.Range("N113").Value = $null
$ws.Range("H122").Value = 6342.8
$ws.Range("I122").Value = 3919.3125
$ws.Range("J122").Value = 16036.75
$ws.Range("K122").Value = 11757.9375
$ws.Range("L122").Value = 48110.25
$ws.Range("M122").Value = -9307.9375
$ws.Range("N122").Value = -53010.25
$ws.Range("H126").Value = 10001.846
$ws.Range("I126").Value = 14308.3
$ws.Range("J126").Value = 7310.3125
$ws.Range("K126").Value = 42924.89999999999
$ws.Range("L126").Value = 21930.9375
$ws.Range("M126").Value = -40454.89999999999
$ws.Range("N126").Value = -26870.9375

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23618.92
$ws.Range("I40").Value = 29066.666
$ws.Range("J40").Value = 15447.3
$ws.Range("K40").Value = 29066.666
$ws.Range("L40").Value = 15447.3
$ws.Range("M40").Value = -28930.666
$ws.Range("N40").Value = -15719.3
$ws.Range("H55").Value = 2054.3333
$ws.Range("I55").Value = 358
$ws.Range("J55").Value = 4174.75
$ws.Range("K55").Value = 358
$ws.Range("L55").Value = 4174.75
$ws.Range("M55").Value = -185
$ws.Range("N55").Value = -4520.75
$ws.Range("H93").Value = 7435.263
$ws.Range("I93").Value = 8422.929
$ws.Range("J93").Value = 4669.8
$ws.Range("K93").Value = 8422.929
$ws.Range("L93").Value = 4669.8
$ws.Range("M93").Value = -7174.929
$ws.Range("N93").Value = -7165.8
$ws.Range("H122").Value = 4573.6787
$ws.Range("I122").Value = 4142.5557
$ws.Range("K122").Value = 12427.6671
$ws.Range("M122").Value = -9977.667099999999
$ws.Range("H136").Value = 5295.24
$ws.Range("I136").Value = 2730.5625
$ws.Range("K136").Value = 8191.6875
$ws.Range("M136").Value = -5641.6875

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1220.079
$ws.Range("I113").Value = 667.36
$ws.Range("J113").Value = 2283
$ws.Range("K113").Value = 2002.08
$ws.Range("L113").Value = 6849
$ws.Range("M113").Value = 167.9200000000001
$ws.Range("N113").Value = -11189
$ws.Range("H122").Value = 4565.1284
$ws.Range("I122").Value = 1695.7693
$ws.Range("J122").Value = 10303.846
$ws.Range("K122").Value = 5087.3079
$ws.Range("L122").Value = 30911.538
$ws.Range("M122").Value = -2637.3079
$ws.Range("N122").Value = -35811.538
$ws.Range("H126").Value = 28711.812
$ws.Range("J126").Value = 6215.3335
$ws.Range("L126").Value = 18646.0005
$ws.Range("N126").Value = -23586.0005
$ws.Range("H136").Value = 418816.78
$ws.Range("I136").Value = 702085.5600000001
$ws.Range("K136").Value = 2106256.68
$ws.Range("M136").Value = -2103706.68
